$d = $word.ActiveDocument

# Fix the title text: "DOCUMENTO SRSS" -> "DOCUMENTO SRS"
$d.Content.Find.Execute("DOCUMENTO SRSS", $true, $false, $false, $false, $false,
                         $true, 1, $false, "DOCUMENTO SRS", 2)

# Underline the paragraph (applies to the paragraph mark, i.e. w:pPr/w:rPr)
$para = $d.Paragraphs.Item(1)
$para.Range.Font.Underline = 1

# The statement above also stamps the explicit underline onto the run that
# holds the visible text. Remove that run-level formatting again by
# deleting the text and re-inserting it, so only the paragraph mark keeps
# the underline (matching the intended "underline the title" edit).
$textRange = $d.Range($para.Range.Start, $para.Range.End - 1)
$title = $textRange.Text
$textRange.Delete()
$insertionPoint = $d.Range($para.Range.Start, $para.Range.Start)
$insertionPoint.InsertBefore($title)
